# Add the ngHO calibration maps Emap/Lmap
#
# - uhtr_side_c26 / uhtr_side_c27: column E (rows 2-7) filled in with the
#   "4-1-xx" calibration map labels (previously a placeholder string).
# - det_side_c26: column E (rows 2-7) filled in with the "HO2M.. RM5" Lmap
#   labels (previously a placeholder "XX" string).
# - det_side_c27: column E (rows 2-7) filled in with the "HO2P.. RM5" Emap
#   labels (previously the same placeholder).
# - Various selection/active-sheet bookkeeping updates left behind by the
#   author while navigating the workbook.

$wb = $excel.ActiveWorkbook

# --- uhtr_side_c26: fill in the E column map labels ---------------------
$wsUhtrC26 = $wb.Worksheets.Item("uhtr_side_c26")
$wsUhtrC26.Range("E2").Value = "4-1-12"
$wsUhtrC26.Range("E3").Value = "4-1-11"
$wsUhtrC26.Range("E4").Value = "4-1-10"
$wsUhtrC26.Range("E5").Value = "4-1-09"
$wsUhtrC26.Range("E6").Value = "4-1-08"
$wsUhtrC26.Range("E7").Value = "4-1-07"

# --- uhtr_side_c27: same E column map labels -----------------------------
$wsUhtrC27 = $wb.Worksheets.Item("uhtr_side_c27")
$wsUhtrC27.Range("E2").Value = "4-1-12"
$wsUhtrC27.Range("E3").Value = "4-1-11"
$wsUhtrC27.Range("E4").Value = "4-1-10"
$wsUhtrC27.Range("E5").Value = "4-1-09"
$wsUhtrC27.Range("E6").Value = "4-1-08"
$wsUhtrC27.Range("E7").Value = "4-1-07"

# --- det_side_c26: fill in the Lmap (HO2M.. RM5) labels ------------------
$wsDetC26 = $wb.Worksheets.Item("det_side_c26")
$wsDetC26.Range("E2").Value = "HO2M12 RM5"
$wsDetC26.Range("E3").Value = "HO2M10 RM5"
$wsDetC26.Range("E4").Value = "HO2M08 RM5"
$wsDetC26.Range("E5").Value = "HO2M06 RM5"
$wsDetC26.Range("E6").Value = "HO2M04 RM5"
$wsDetC26.Range("E7").Value = "HO2M02 RM5"

# --- det_side_c27: fill in the Emap (HO2P.. RM5) labels ------------------
$wsDetC27 = $wb.Worksheets.Item("det_side_c27")
$wsDetC27.Range("E2").Value = "HO2P12 RM5"
$wsDetC27.Range("E3").Value = "HO2P10 RM5"
$wsDetC27.Range("E4").Value = "HO2P08 RM5"
$wsDetC27.Range("E5").Value = "HO2P06 RM5"
$wsDetC27.Range("E6").Value = "HO2P04 RM5"
$wsDetC27.Range("E7").Value = "HO2P02 RM5"

# --- selection / active-sheet bookkeeping --------------------------------
# Leave uhtr_side_c23 alone (no selection change in the diff).

# uhtr_side_c26: cursor moved to M19
[void]$wsUhtrC26.Range("M19").Select()

# uhtr_side_c27: cursor moved to the newly filled E2:E7 block
[void]$wsUhtrC27.Range("E2:E7").Select()

# uhtr_side_c33: cursor moved to the E2:E7 block as well
$wsUhtrC33 = $wb.Worksheets.Item("uhtr_side_c33")
[void]$wsUhtrC33.Range("E2:E7").Select()

# det_side_c27 becomes the active / tab-selected sheet, cursor on E8 --
# select it last so it ends up as the workbook's active tab.
[void]$wsDetC27.Range("E8").Select()

Write-Output "applied ngHO calibration map edits"
